$d = $word.ActiveDocument

# Replace the text of the first 21 paragraphs in place (1:1 swap),
# setting Range.Text directly so no AutoCorrect/"smart quote" substitution occurs.
$d.Paragraphs(1).Range.Text = "Module 09: Tissues and the Animal Body — Study Questions"
$d.Paragraphs(2).Range.Text = "What is homeostasis? Give an example of a condition the human body must keep perfectly balanced."
$d.Paragraphs(3).Range.Text = "How does a negative feedback loop differ from a positive feedback loop? Provide a biological example of each."
$d.Paragraphs(4).Range.Text = "Why is osmoregulation critical for human survival? What would happen to our cells if this failed?"
$d.Paragraphs(5).Range.Text = "Which organ system plays the largest role in maintaining water and salt balance (osmoregulation)?"
$d.Paragraphs(6).Range.Text = "What is the fundamental difference between mechanical digestion and chemical digestion?"
$d.Paragraphs(7).Range.Text = "Trace the path of food from the mouth to the large intestine, naming the primary function of each organ it passes through."
$d.Paragraphs(8).Range.Text = "Why is the small intestine considered the most important organ for digestion and absorption?"
$d.Paragraphs(9).Range.Text = "The liver, pancreas, and gallbladder are `"accessory organs.`" What does this mean, and what specific fluids do they contribute?"
$d.Paragraphs(10).Range.Text = "What are the three main types of blood vessels, and how does their structure match their function?"
$d.Paragraphs(11).Range.Text = "Trace the pathway of a red blood cell starting from the right atrium of the heart, to the lungs, and back to the left side of the heart."
$d.Paragraphs(12).Range.Text = "Why is the left ventricle of the heart thicker and more muscular than the right ventricle?"
$d.Paragraphs(13).Range.Text = "Describe the process of gas exchange in the alveoli. What gas moves into the blood, and what gas moves out?"
$d.Paragraphs(14).Range.Text = "What is a hormone, and how does the endocrine system use hormones differently than how the nervous system uses electrical signals?"
$d.Paragraphs(15).Range.Text = "Choose two endocrine glands (e.g., pancreas, thyroid) and describe their primary function or the specific hormone they produce."
$d.Paragraphs(16).Range.Text = "How does the pancreas use hormones to regulate blood sugar levels?"
$d.Paragraphs(17).Range.Text = "Which part of the brain acts as the primary link between the nervous system and the endocrine system?"
$d.Paragraphs(18).Range.Text = "Besides providing structural support and movement, what are two other critical functions of the human skeleton?"
$d.Paragraphs(19).Range.Text = "What is the difference between the axial skeleton and the appendicular skeleton?"
$d.Paragraphs(20).Range.Text = "Compare and contrast the three types of muscle tissue (skeletal, smooth, cardiac) in terms of their location and whether their movement is voluntary or involuntary."
$d.Paragraphs(21).Range.Text = "Using an analogy of a lever, explain how skeletal muscles and bones interact to create movement."

# Append 4 new question paragraphs at the end of the document
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Text = "What are the two main anatomical divisions of the nervous system, and what structures belong to each?"
$d.Paragraphs($count).Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Text = "Draw or describe the basic structure of a neuron, including the dendrites, cell body, and axon."
$d.Paragraphs($count).Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Text = "What is the difference between a sensory neuron and a motor neuron?"
$d.Paragraphs($count).Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Text = "Challenge Question: You touch a hot stove and immediately pull your hand away. Describe how the nervous, muscular, and skeletal systems worked together almost instantly to protect you from a severe burn."

Write-Output "done"